# Update cryptocurrency price/volume data (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "60.734.11"
$ws.Range("E2").Value = "  +3.14%  "

# Row 3
$ws.Range("D3").Value = "2.685.13"
$ws.Range("E3").Value = "  +1.98%  "

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").Value = "'522.19"
$ws.Range("E5").Value = "  +1.68%  "

# Row 6
$ws.Range("D6").Value = "'147.33"
$ws.Range("E6").Value = "  +1.89%  "

# Row 7
$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "  +0.10%  "

# Row 8
$ws.Range("D8").Value = "'0.580"
$ws.Range("E8").Value = "  +2.21%  "

# Row 9
$ws.Range("D9").Value = "2.704.49"
$ws.Range("E9").Value = "  +1.67%  "

# Row 10
$ws.Range("D10").Value = "'6.45"
$ws.Range("E10").Value = "  +0.60%  "

# Row 11
$ws.Range("E11").Value = "  +0.92%  "

# Row 12
$ws.Range("D12").Value = "'0.342"
$ws.Range("E12").Value = "  +1.18%  "

# Row 13
$ws.Range("E13").Value = "  +1.47%  "

# Row 14
$ws.Range("D14").Value = "3.158.74"
$ws.Range("E14").Value = "  +2.13%  "

# Row 15
$ws.Range("D15").Value = "60.762.58"
$ws.Range("E15").Value = "  +3.20%  "

# Row 16
$ws.Range("D16").Value = "'21.45"
$ws.Range("E16").Value = "  +1.65%  "

# Row 17
$ws.Range("D17").Value = "2.762.39"
$ws.Range("E17").Value = "  +3.16%  "

# Row 18
$ws.Range("E18").Value = "  +1.45%  "

# Row 19
$ws.Range("D19").Value = "'353.58"
$ws.Range("E19").Value = "  +2.50%  "

# Row 20
$ws.Range("D20").Value = "'4.58"
$ws.Range("E20").Value = "  +0.39%  "

# Row 21
$ws.Range("D21").Value = "'10.58"
$ws.Range("E21").Value = "  +1.84%  "

# Row 22
$ws.Range("D22").Value = "'6.37"
$ws.Range("E22").Value = "  +4.12%  "

# Row 23
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.27%  "

# Row 24
$ws.Range("D24").Value = "'63.08"
$ws.Range("E24").Value = "  +2.70%  "

# Row 25
$ws.Range("D25").Value = "'0.424"
$ws.Range("E25").Value = "  +0.43%  "

# Row 26
$ws.Range("D26").Value = "'0.170"
$ws.Range("E26").Value = "  +4.77%  "

# Row 27
$ws.Range("D27").Value = "'0.994"
$ws.Range("E27").Value = "  -0.10%  "

# Row 28
$ws.Range("D28").Value = "0.0₃0821"
$ws.Range("E28").Value = "  +1.46%  "

# Row 29
$ws.Range("E29").Value = "  +2.50%  "

# Row 30
$ws.Range("D30").Value = "'6.91"
$ws.Range("E30").Value = "  +6.76%  "

# Row 31
$ws.Range("E31").Value = "  +0.11%  "

# Row 32
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'1.60"
$ws.Range("E32").Value = "  +1.23%  "

# Row 33
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "'19.14"
$ws.Range("E33").Value = "  +1.04%  "

# Row 34
$ws.Range("D34").Value = "'149.48"
$ws.Range("E34").Value = "  -0.39%  "

# Row 35
$ws.Range("D35").Value = "'4.34"
$ws.Range("E35").Value = "  +7.57%  "

# Row 36
$ws.Range("B36").Value = "SuiNetwork"
$ws.Range("C36").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D36").Value = "'0.953"
$ws.Range("E36").Value = "  -7.95%  "

# Row 37
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'1.24"
$ws.Range("E37").Value = "  +6.76%  "

# Row 38
$ws.Range("E38").Value = "  +10.90%  "

# Row 39
$ws.Range("D39").Value = "'0.880"
$ws.Range("E39").Value = "  +2.48%  "

# Row 40
$ws.Range("D40").Value = "'36.74"
$ws.Range("E40").Value = "  +0.50%  "

# Row 41
$ws.Range("E41").Value = "  +0.65%  "

# Row 42
$ws.Range("D42").Value = "'285.56"
$ws.Range("E42").Value = "  +1.69%  "

# Row 43
$ws.Range("D43").Value = "'20.09"
$ws.Range("E43").Value = "  +2.95%  "

# Row 44
$ws.Range("D44").Value = "'0.0994"
$ws.Range("E44").Value = "  +1.03%  "

# Row 45
$ws.Range("D45").Value = "'0.613"
$ws.Range("E45").Value = "  -0.14%  "

# Row 46
$ws.Range("D46").Value = "'0.996"
$ws.Range("E46").Value = "  +0.19%  "

# Row 47
$ws.Range("D47").Value = "2.142.45"
$ws.Range("E47").Value = "  +8.06%  "

# Row 48
$ws.Range("D48").Value = "'0.0543"
$ws.Range("E48").Value = "  +0.78%  "

# Row 49
$ws.Range("D49").Value = "'4.89"
$ws.Range("E49").Value = "  +4.23%  "

# Row 50
$ws.Range("E50").Value = "  +2.90%  "

# Row 51
$ws.Range("D51").Value = "'10.47"
$ws.Range("E51").Value = "  +2.09%  "

